{"js": "// \"Remove Internet of Things\" \u2014 two small wording tweaks in the Experience\n// section of the resume:\n//   1. \"atSistemas for IBM Software\" -> \"atSistemas for IBM\"\n//   2. \"a Smart Cities/Internet of Things project\" -> \"a Smart Cities\"\n\nconst body = context.document.body;\n\n// 1) \"atSistemas for IBM Software\" -> \"atSistemas for IBM\"\nconst titleResults = body.search(\"atSistemas for IBM Software\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"atSistemas for IBM\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"a Smart Cities/Internet of Things project\" -> \"a Smart Cities\"\nconst projectResults = body.search(\"a Smart Cities/Internet of Things project\", { matchCase: true });\nprojectResults.load(\"items\");\nawait context.sync();\n\nif (projectResults.items.length > 0) {\n  projectResults.items[0].insertText(\"a Smart Cities\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Remove Internet of Things\" \u2014 two small wording tweaks in the Experience\n# section of the resume:\n#   1. \"atSistemas for IBM Software\" -> \"atSistemas for IBM\"\n#   2. \"a Smart Cities/Internet of Things project\" -> \"a Smart Cities\"\n\n$d = $word.ActiveDocument\n\n# 1) \"atSistemas for IBM Software\" -> \"atSistemas for IBM\"\n$r1 = $d.Content\n$r1.Find.Text = \"atSistemas for IBM Software\"\nif ($r1.Find.Execute()) {\n    $r1.Text = \"atSistemas for IBM\"\n}\n\n# 2) \"a Smart Cities/Internet of Things project\" -> \"a Smart Cities\"\n$r2 = $d.Content\n$r2.Find.Text = \"a Smart Cities/Internet of Things project\"\nif ($r2.Find.Execute()) {\n    $r2.Text = \"a Smart Cities\"\n}\n"}
